# TC17_Canine_Filter_Breed-Chihuahua.xlsx - "Fixed ICDC breed all testcases"
#
# The StatQuery column (C) on the "startup" sheet previously held an old,
# malformed Cypher aggregate query (shared across rows 2-4: Cases/Samples/
# Files). It is replaced here with a corrected query that counts Programs,
# Studies, Cases, Samples, Case Files and Study Files. Columns B (the
# per-tab query) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Chihuahua']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Window/view tweak that accompanied the fix: zoomed in a bit and moved the
# active selection down to the FilesTab row.
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("B4").Select()
